$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 565.8182
$ws.Range("I38").Value = 104
$ws.Range("J38").Value = 1120
$ws.Range("K38").Value = 312
$ws.Range("L38").Value = 3360
$ws.Range("M38").Value = 60
$ws.Range("N38").Value = -4104

$ws.Range("H41").Value = 2675
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()

$ws.Range("H42").Value = 93.5
$ws.Range("I42").Value = 93.5
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 280.5
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = -50.5
$ws.Range("N42").ClearContents()

$ws.Range("H43").Value = 291668930
$ws.Range("J43").Value = 55558560
$ws.Range("L43").Value = 55558560
$ws.Range("N43").Value = -55558698

$ws.Range("H76").Value = 127837.875
$ws.Range("I76").Value = 145600.42
$ws.Range("J76").Value = 3500
$ws.Range("K76").Value = 145600.42
$ws.Range("L76").Value = 3500
$ws.Range("M76").Value = -145285.42
$ws.Range("N76").Value = -4130

$ws.Range("H79").Value = 127837.875
$ws.Range("I79").Value = 145600.42
$ws.Range("J79").Value = 3500
$ws.Range("K79").Value = 145600.42
$ws.Range("L79").Value = 3500
$ws.Range("M79").Value = -144508.42
$ws.Range("N79").Value = -5684

$ws.Range("H95").Value = 37257
$ws.Range("J95").Value = 37257
$ws.Range("L95").Value = 37257
$ws.Range("N95").Value = -42749

$ws.Range("H97").Value = 1716.5555
$ws.Range("I97").Value = 224.5
$ws.Range("J97").Value = 2142.8572
$ws.Range("K97").Value = 673.5
$ws.Range("L97").Value = 6428.571599999999
$ws.Range("M97").Value = -177.5
$ws.Range("N97").Value = -7420.571599999999

$ws.Range("H100").Value = 1253.3077
$ws.Range("I100").Value = 1349.125
$ws.Range("J100").Value = 1100
$ws.Range("K100").Value = 1349.125
$ws.Range("L100").Value = 1100
$ws.Range("M100").Value = -808.125
$ws.Range("N100").Value = -2182

$ws.Range("H113").Value = 3109.1765
$ws.Range("I113").Value = 2936
$ws.Range("J113").Value = 3263.111
$ws.Range("K113").Value = 2936
$ws.Range("L113").Value = 3263.111
$ws.Range("M113").Value = 318
$ws.Range("N113").Value = -9771.111000000001

$ws.Range("H137").Value = 1154.6724
$ws.Range("I137").Value = 1015.3555
$ws.Range("J137").Value = 1636.9231
$ws.Range("K137").Value = 3046.0665
$ws.Range("L137").Value = 4910.7693
$ws.Range("M137").Value = -496.0664999999999
$ws.Range("N137").Value = -10010.7693

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1481.6
$ws.Range("I2").Value = 1301.5714
$ws.Range("J2").Value = 1639.125
$ws.Range("K2").Value = 1301.5714
$ws.Range("L2").Value = 1639.125
$ws.Range("M2").Value = -1188.5714
$ws.Range("N2").Value = -1865.125

$ws.Range("H37").Value = 12025.308
$ws.Range("I37").Value = 8237.5
$ws.Range("J37").Value = 13708.777
$ws.Range("K37").Value = 8237.5
$ws.Range("L37").Value = 13708.777
$ws.Range("M37").Value = -7964.5
$ws.Range("N37").Value = -14254.777

$ws.Range("H97").Value = 2104.2
$ws.Range("I97").Value = 2336.6667
$ws.Range("J97").Value = 1755.5
$ws.Range("K97").Value = 2336.6667
$ws.Range("L97").Value = 1755.5
$ws.Range("M97").Value = -1840.6667
$ws.Range("N97").Value = -2747.5

$ws.Range("H98").Value = 19000
$ws.Range("J98").Value = 19000
$ws.Range("L98").Value = 19000
$ws.Range("N98").Value = -24990

$ws.Range("H116").Value = 1481.6
$ws.Range("I116").Value = 1301.5714
$ws.Range("J116").Value = 1639.125
$ws.Range("K116").Value = 1301.5714
$ws.Range("L116").Value = 1639.125
$ws.Range("M116").Value = 992.4286
$ws.Range("N116").Value = -6227.125

$ws.Range("H132").Value = 8793.532999999999
$ws.Range("I132").Value = 9908.75
$ws.Range("J132").Value = 4332.6665
$ws.Range("K132").Value = 29726.25
$ws.Range("L132").Value = 12997.9995
$ws.Range("M132").Value = -27196.25
$ws.Range("N132").Value = -18057.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1481.6
$ws.Range("I3").Value = 1301.5714
$ws.Range("J3").Value = 1639.125
$ws.Range("K3").Value = 1301.5714
$ws.Range("L3").Value = 1639.125
$ws.Range("M3").Value = -1187.5714
$ws.Range("N3").Value = -1867.125

$ws.Range("H94").Value = 3500
$ws.Range("I94").Value = 2000
$ws.Range("J94").Value = 5000
$ws.Range("K94").Value = 2000
$ws.Range("L94").Value = 5000
$ws.Range("M94").Value = -1549
$ws.Range("N94").Value = -5902

$ws.Range("H99").Value = 720.1111
$ws.Range("I99").Value = 710
$ws.Range("J99").Value = 755.5
$ws.Range("K99").Value = 710
$ws.Range("L99").Value = 755.5
$ws.Range("M99").Value = 788
$ws.Range("N99").Value = -3751.5

$ws.Range("H107").Value = 1245
$ws.Range("I107").Value = 1250
$ws.Range("J107").Value = 1240
$ws.Range("K107").Value = 1250
$ws.Range("L107").Value = 1240
$ws.Range("M107").Value = 670
$ws.Range("N107").Value = -5080

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H96").Value = 29399.5
$ws.Range("J96").Value = 29399.5
$ws.Range("L96").Value = 29399.5
$ws.Range("N96").Value = -34891.5

$ws.Range("H132").Value = 3516.9583
$ws.Range("I132").Value = 2861.6667
$ws.Range("J132").Value = 5482.8335
$ws.Range("K132").Value = 8585.000100000001
$ws.Range("L132").Value = 16448.5005
$ws.Range("M132").Value = -6055.000100000001
$ws.Range("N132").Value = -21508.5005

$ws.Range("H134").Value = 1016.13635
$ws.Range("I134").Value = 972.4865
$ws.Range("J134").Value = 1246.8572
$ws.Range("K134").Value = 2917.4595
$ws.Range("L134").Value = 3740.5716
$ws.Range("M134").Value = -382.4594999999999
$ws.Range("N134").Value = -8810.571599999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 1109.75
$ws.Range("I109").Value = 649.6667
$ws.Range("J109").Value = 2490
$ws.Range("K109").Value = 1949.0001
$ws.Range("L109").Value = 7470
$ws.Range("M109").Value = -909.0001
$ws.Range("N109").Value = -9550

$ws.Range("H113").Value = 554.0417
$ws.Range("I113").Value = 524.28
$ws.Range("K113").Value = 1572.84
$ws.Range("M113").Value = 597.1600000000001

$ws.Range("H131").Value = 6542390
$ws.Range("I131").Value = 9529
$ws.Range("J131").Value = 18519302
$ws.Range("K131").Value = 28587
$ws.Range("L131").Value = 55557906
$ws.Range("M131").Value = -23547
$ws.Range("N131").Value = -55567986

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 21128.6
$ws.Range("J95").Value = 21128.6
$ws.Range("L95").Value = 21128.6
$ws.Range("N95").Value = -26620.6

$ws.Range("H107").Value = 766.2727
$ws.Range("I107").Value = 730
$ws.Range("J107").Value = 818.6667
$ws.Range("K107").Value = 730
$ws.Range("L107").Value = 818.6667
$ws.Range("M107").Value = 1190
$ws.Range("N107").Value = -4658.6667

$ws.Range("H132").Value = 54978.973
$ws.Range("I132").Value = 66354
$ws.Range("J132").Value = 4603.857
$ws.Range("K132").Value = 199062
$ws.Range("L132").Value = 13811.571
$ws.Range("M132").Value = -196532
$ws.Range("N132").Value = -18871.571

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 460.625
$ws.Range("I55").Value = 475
$ws.Range("J55").Value = 446.25
$ws.Range("K55").Value = 475
$ws.Range("L55").Value = 446.25
$ws.Range("M55").Value = -302
$ws.Range("N55").Value = -792.25

$ws.Range("H61").Value = 948.2778
$ws.Range("I61").Value = 948.2778
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 948.2778
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -746.2778
$ws.Range("N61").ClearContents()

$ws.Range("H93").Value = 1554.5714
$ws.Range("I93").Value = 1617
$ws.Range("J93").Value = 1398.5
$ws.Range("K93").Value = 1617
$ws.Range("L93").Value = 1398.5
$ws.Range("M93").Value = -369
$ws.Range("N93").Value = -3894.5

$ws.Range("H113").Value = 948.2778
$ws.Range("I113").Value = 948.2778
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 948.2778
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1221.7222
$ws.Range("N113").ClearContents()

$ws.Range("H132").Value = 1698.8975
$ws.Range("I132").Value = 1024.5161
$ws.Range("K132").Value = 3073.5483
$ws.Range("M132").Value = -543.5483000000004

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 54800
$ws.Range("J94").Value = 54800
$ws.Range("L94").Value = 54800
$ws.Range("N94").Value = -56602

$ws.Range("H100").Value = 790
$ws.Range("I100").Value = 740
$ws.Range("J100").Value = 860
$ws.Range("K100").Value = 1480
$ws.Range("L100").Value = 1720
$ws.Range("M100").Value = -939
$ws.Range("N100").Value = -2802

$ws.Range("H136").Value = 3441.125
$ws.Range("I136").Value = 4006.3242
$ws.Range("K136").Value = 12018.9726
$ws.Range("M136").Value = -9468.972600000001
